$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "2022-Q4" sheet.
#    It carries the same fund roster as "2022-Q3" (same fund codes/names),
#    just with refreshed position numbers, so clone the "2022-Q3" sheet
#    (preserves headers, styles, text-vs-number cell typing) and place the
#    clone immediately in front of it, then touch up the changed figures.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Row 2 - 090019 大成景恒混合A
$q4.Cells.Item(2,4).Formula = "'1.18"
$q4.Cells.Item(2,5).Formula = "'93.72"
$q4.Cells.Item(2,6).Formula = "'1.62"
$q4.Cells.Item(2,7).Formula = "'0.0191"
$q4.Cells.Item(2,8).Value = 8

# Row 3 - 006038 大成景恒混合C
$q4.Cells.Item(3,4).Formula = "'0.89"
$q4.Cells.Item(3,5).Formula = "'93.72"
$q4.Cells.Item(3,6).Formula = "'1.62"
$q4.Cells.Item(3,7).Formula = "'0.0144"
$q4.Cells.Item(3,8).Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a row for 2022-Q4 at the top
#    of the table and push the existing quarters down by one row.
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Keep a correctly-styled template (column A uses a bold/bordered/centered
# style) so newly written cells in column A match the existing look without
# fabricating a brand new style entry.
$zj.Cells.Item(3,1).Copy()

# Walk from the bottom up so source rows are read before being overwritten.

# Row 6 <- was row 5 (2021-Q2)
$zj.Cells.Item(6,1).PasteSpecial(-4122)
$zj.Cells.Item(6,1).Value = 4
$zj.Cells.Item(6,2).Value = "2021-Q2"
$zj.Cells.Item(6,3).Value = 22
$zj.Cells.Item(6,4).Value = 3.74

# Row 5 <- was row 4 (2021-Q3)
$zj.Cells.Item(5,1).PasteSpecial(-4122)
$zj.Cells.Item(5,1).Value = 3
$zj.Cells.Item(5,2).Value = "2021-Q3"
$zj.Cells.Item(5,3).Value = 9
$zj.Cells.Item(5,4).Value = 2.62

# Row 4 <- was row 3 (2021-Q4)
$zj.Cells.Item(4,1).Value = 2
$zj.Cells.Item(4,2).Value = "2021-Q4"
$zj.Cells.Item(4,3).Value = 10
$zj.Cells.Item(4,4).Value = 1.52

# Row 3 <- was row 2 (2022-Q3)
$zj.Cells.Item(3,1).Value = 1
$zj.Cells.Item(3,2).Value = "2022-Q3"
$zj.Cells.Item(3,3).Value = 2
$zj.Cells.Item(3,4).Value = 0.03

# Row 2 <- new 2022-Q4 entry (holdings unchanged versus Q3: 2 funds, 0.03亿)
$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = "2022-Q4"
$zj.Cells.Item(2,3).Value = 2
$zj.Cells.Item(2,4).Value = 0.03
